$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 1098.7
$ws.Range("I82").Value = 263.7143
$ws.Range("J82").Value = 3047
$ws.Range("K82").Value = 791.1428999999999
$ws.Range("L82").Value = 9141
$ws.Range("M82").Value = -385.1428999999999
$ws.Range("N82").Value = -9953
$ws.Range("H85").Value = 1098.7
$ws.Range("I85").Value = 263.7143
$ws.Range("J85").Value = 3047
$ws.Range("K85").Value = 791.1428999999999
$ws.Range("L85").Value = 9141
$ws.Range("M85").Value = 612.8571000000001
$ws.Range("N85").Value = -11949
$ws.Range("H112").Value = 1112.18
$ws.Range("J112").Value = 1135.699
$ws.Range("L112").Value = 3407.097
$ws.Range("N112").Value = -5623.097
$ws.Range("H113").Value = 2017.1904
$ws.Range("I113").Value = 1902.2727
$ws.Range("J113").Value = 2143.6
$ws.Range("K113").Value = 1902.2727
$ws.Range("L113").Value = 2143.6
$ws.Range("M113").Value = 1351.7273
$ws.Range("N113").Value = -8651.6
$ws.Range("H132").Value = 195234.97
$ws.Range("I132").Value = 2958.4285
$ws.Range("J132").Value = 1002796.4
$ws.Range("K132").Value = 8875.2855
$ws.Range("L132").Value = 3008389.2
$ws.Range("M132").Value = -6345.2855
$ws.Range("N132").Value = -3013449.2
$ws.Range("H138").Value = 2286.4
$ws.Range("I138").Value = 994.7143
$ws.Range("J138").Value = 2629.7595
$ws.Range("K138").Value = 2984.1429
$ws.Range("L138").Value = 7889.2785
$ws.Range("M138").Value = 2155.8571
$ws.Range("N138").Value = -18169.2785

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 41801.08
$ws.Range("I74").Value = 49000.383
$ws.Range("J74").Value = 4004.75
$ws.Range("K74").Value = 49000.383
$ws.Range("L74").Value = 4004.75
$ws.Range("M74").Value = -48126.383
$ws.Range("N74").Value = -5752.75
$ws.Range("H77").Value = 41801.08
$ws.Range("I77").Value = 49000.383
$ws.Range("J77").Value = 4004.75
$ws.Range("K77").Value = 245001.915
$ws.Range("L77").Value = 20023.75
$ws.Range("M77").Value = -240633.915
$ws.Range("N77").Value = -28759.75
$ws.Range("H132").Value = 1472473.2
$ws.Range("I132").Value = 1673086.6
$ws.Range("J132").Value = 916227
$ws.Range("K132").Value = 5019259.800000001
$ws.Range("L132").Value = 2748681
$ws.Range("M132").Value = -5016729.800000001
$ws.Range("N132").Value = -2753741

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 293164.66
$ws.Range("I86").Value = 1511.6154
$ws.Range("J86").Value = 637845.5600000001
$ws.Range("K86").Value = 1511.6154
$ws.Range("L86").Value = 637845.5600000001
$ws.Range("M86").Value = -388.6153999999999
$ws.Range("N86").Value = -640091.5600000001
$ws.Range("H89").Value = 293164.66
$ws.Range("I89").Value = 1511.6154
$ws.Range("J89").Value = 637845.5600000001
$ws.Range("K89").Value = 7558.076999999999
$ws.Range("L89").Value = 3189227.8
$ws.Range("M89").Value = -1942.076999999999
$ws.Range("N89").Value = -3200459.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10348
$ws.Range("H25").Value = 29919.625
$ws.Range("I25").Value = 1500
$ws.Range("J25").Value = 39392.832
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 39392.832
$ws.Range("N25").Value = -39740.832
$ws.Range("H41").Value = 10200
$ws.Range("J41").Value = 19900
$ws.Range("L41").Value = 19900
$ws.Range("H50").Value = 40900
$ws.Range("I50").Value = 40900
$ws.Range("K50").Value = 40900
$ws.Range("M50").Value = -40275
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("H59").Value = 11444.444
$ws.Range("J59").Value = 11444.444
$ws.Range("L59").Value = 11444.444
$ws.Range("N59").Value = -13734.444
$ws.Range("H60").Value = 7057.2
$ws.Range("I60").Value = 4321.5
$ws.Range("J60").Value = 18000
$ws.Range("K60").Value = 4321.5
$ws.Range("L60").Value = 18000
$ws.Range("M60").Value = -3810.5
$ws.Range("N60").Value = -19022
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("H68").Value = 29666.666
$ws.Range("J68").Value = 29666.666
$ws.Range("L68").Value = 29666.666
$ws.Range("N68").Value = -31164.666
$ws.Range("H71").Value = 29666.666
$ws.Range("J71").Value = 29666.666
$ws.Range("L71").Value = 88999.99800000001
$ws.Range("N71").Value = -96487.99800000001
$ws.Range("H74").Value = 12211.111
$ws.Range("J74").Value = 12211.111
$ws.Range("L74").Value = 12211.111
$ws.Range("N74").Value = -13959.111
$ws.Range("H77").Value = 12211.111
$ws.Range("J77").Value = 12211.111
$ws.Range("L77").Value = 36633.333
$ws.Range("N77").Value = -45369.333
$ws.Range("M51").ClearContents()
$ws.Range("M61").ClearContents()
$ws.Range("M25").Value = -1326
$ws.Range("N41").Value = -20756

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1242.3334
$ws.Range("I107").Value = 821.1177
$ws.Range("J107").Value = 1619.2106
$ws.Range("K107").Value = 2463.3531
$ws.Range("L107").Value = 4857.6318
$ws.Range("M107").Value = -543.3531000000003
$ws.Range("N107").Value = -8697.631799999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7291.1
$ws.Range("I80").Value = 2578.5715
$ws.Range("J80").Value = 9828.615
$ws.Range("K80").Value = 2578.5715
$ws.Range("L80").Value = 9828.615
$ws.Range("M80").Value = -1580.5715
$ws.Range("N80").Value = -11824.615
$ws.Range("H83").Value = 7291.1
$ws.Range("I83").Value = 2578.5715
$ws.Range("J83").Value = 9828.615
$ws.Range("K83").Value = 12892.8575
$ws.Range("L83").Value = 49143.075
$ws.Range("M83").Value = -7900.8575
$ws.Range("N83").Value = -59127.075

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1645.9688
$ws.Range("I68").Value = 1515.2106
$ws.Range("J68").Value = 1837.0769
$ws.Range("K68").Value = 1515.2106
$ws.Range("L68").Value = 1837.0769
$ws.Range("M68").Value = -766.2106000000001
$ws.Range("N68").Value = -3335.0769
$ws.Range("H71").Value = 1645.9688
$ws.Range("I71").Value = 1515.2106
$ws.Range("J71").Value = 1837.0769
$ws.Range("K71").Value = 7576.053000000001
$ws.Range("L71").Value = 9185.3845
$ws.Range("M71").Value = -3832.053000000001
$ws.Range("N71").Value = -16673.3845

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4286.515
$ws.Range("I122").Value = 3347.4
$ws.Range("J122").Value = 5731.3076
$ws.Range("K122").Value = 10042.2
$ws.Range("L122").Value = 17193.9228
$ws.Range("M122").Value = -7592.200000000001
$ws.Range("N122").Value = -22093.9228
